$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G3").Value = 2.6
$ws.Range("I3").Value = 3.65
$ws.Range("L3").Value = 1.34
$ws.Range("T4").Value = 1.75
$ws.Range("U4").Value = 2.08
$ws.Range("F5").Value = 1.27
$ws.Range("P5").Value = 2.66
$ws.Range("Q5").Value = 1.53
$ws.Range("R5").Value = 1.66
$ws.Range("S5").Value = 2.32
$ws.Range("F6").Value = 2.44
$ws.Range("G6").Value = 2.5
$ws.Range("I6").Value = 3.45
$ws.Range("J6").Value = 3.3
$ws.Range("S6").Value = 4.3
$ws.Range("W6").Value = 1.67
$ws.Range("AN6").Value = 26
$ws.Range("I7").Value = 8.6
$ws.Range("L7").Value = 1.23
$ws.Range("Q7").Value = 1.51
$ws.Range("R7").Value = 1.62
$ws.Range("S7").Value = 2.1
$ws.Range("X7").Value = 29
$ws.Range("Y7").Value = 32
$ws.Range("AE7").Value = 95
$ws.Range("AM7").Value = 100
$ws.Range("F8").Value = 1.39
$ws.Range("G8").Value = 1.48
$ws.Range("H8").Value = 8.199999999999999
$ws.Range("I8").Value = 11.5
$ws.Range("J8").Value = 4.8
$ws.Range("K8").Value = 6.4
$ws.Range("P8").Value = 2.28
$ws.Range("Q8").Value = 1.63
$ws.Range("S8").Value = 2.58
$ws.Range("W8").Value = 3.05
$ws.Range("X8").Value = 27
$ws.Range("AA8").Value = 350
$ws.Range("AB8").Value = 12
$ws.Range("AC8").Value = 14.5
$ws.Range("AF8").Value = 11.5
$ws.Range("AG8").Value = 13
$ws.Range("AJ8").Value = 15
$ws.Range("AM8").Value = 150
$ws.Range("AN8").Value = 7
$ws.Range("F9").Value = 4.7
$ws.Range("I9").Value = 1.91
$ws.Range("Q9").Value = 1.94
$ws.Range("T9").Value = 1.77
$ws.Range("V9").Value = 2.08
$ws.Range("Q10").Value = 1.43
$ws.Range("V10").Value = 3.3
$ws.Range("G11").Value = 1.55
$ws.Range("H11").Value = 8.199999999999999
$ws.Range("I11").Value = 9.199999999999999
$ws.Range("O11").Value = 1.45
$ws.Range("Q11").Value = 2.3
$ws.Range("R11").Value = 1.25
$ws.Range("T11").Value = 2.42
$ws.Range("U11").Value = 1.61
$ws.Range("F13").Value = 1.48
$ws.Range("G13").Value = 1.5
$ws.Range("K13").Value = 5.1
$ws.Range("S13").Value = 2.86
$ws.Range("V13").Value = 1.13
$ws.Range("W13").Value = 3
$ws.Range("G14").Value = 1.54
$ws.Range("J14").Value = 3.7
$ws.Range("L14").Value = 1.39
$ws.Range("Q14").Value = 2.14
$ws.Range("V14").Value = 1.09
$ws.Range("W14").Value = 2.84
$ws.Range("L15").Value = 1.51
$ws.Range("F16").Value = 1.6
$ws.Range("G16").Value = 1.79
$ws.Range("H16").Value = 6
$ws.Range("I16").Value = 7.6
$ws.Range("K16").Value = 4.7
$ws.Range("N16").Value = 2.98
$ws.Range("O16").Value = 1.38
$ws.Range("P16").Value = 1.68
$ws.Range("Q16").Value = 2.12
$ws.Range("G17").Value = 4.7
$ws.Range("I17").Value = 2.3
$ws.Range("J17").Value = 3.15
$ws.Range("U17").Value = 1.81
$ws.Range("V17").Value = 1.77
$ws.Range("F18").Value = 4.1
$ws.Range("N18").Value = 2.98
$ws.Range("P18").Value = 1.67
$ws.Range("Q18").Value = 2.24
$ws.Range("T18").Value = 1.96
$ws.Range("U18").Value = 1.88
$ws.Range("X18").Value = 11.5
$ws.Range("F19").Value = 2.24
$ws.Range("J19").Value = 3.2
$ws.Range("P19").Value = 1.72
$ws.Range("R19").Value = 1.27
$ws.Range("F20").Value = 1.9
$ws.Range("G20").Value = 1.98
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 1.73
$ws.Range("W20").Value = 2.02
$ws.Range("F21").Value = 1.87
$ws.Range("G21").Value = 1.88
$ws.Range("W21").Value = 2.12
$ws.Range("AG21").Value = 10
$ws.Range("AJ21").Value = 22
$ws.Range("AK21").Value = 16.5
$ws.Range("AL21").Value = 25
$ws.Range("F22").Value = 2.7
$ws.Range("K22").Value = 3.55
$ws.Range("P22").Value = 1.92
$ws.Range("F24").Value = 3.35
$ws.Range("G24").Value = 3.9
$ws.Range("H24").Value = 1.99
$ws.Range("I24").Value = 2.1
$ws.Range("J24").Value = 4
$ws.Range("M24").Value = 1.03
$ws.Range("P24").Value = 2.2
$ws.Range("Q24").Value = 1.64
$ws.Range("S24").Value = 2.44
$ws.Range("T24").Value = 1.62
$ws.Range("U24").Value = 2.24
$ws.Range("V24").Value = 1.9
$ws.Range("W24").Value = 1.34
$ws.Range("AE24").Value = 980
$ws.Range("AF24").Value = 1000
$ws.Range("AG24").Value = 980
$ws.Range("AI24").Value = 980
$ws.Range("AJ24").Value = 85
$ws.Range("AK24").Value = 980
$ws.Range("AL24").Value = 1000
$ws.Range("AM24").Value = 1000
$ws.Range("AN24").Value = 980
